$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.690.32"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "2.534.67"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "565.95"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "146.17"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").Value = "2.534.08"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "26.94"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "2.988.72"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "62.712.93"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "2.528.64"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "11.43"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "333.62"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "64.60"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("D26").Value = "1.59"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "1.48"
$ws.Range("E28").Value = "  +11.32%  "
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +4.65%  "
$ws.Range("D31").Value = "0.0₃0807"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "1.85"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "176.70"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("D35").Value = "401.89"
$ws.Range("E35").Value = "  +7.65%  "
$ws.Range("D36").Value = "0.396"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "18.87"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D39").Value = "4.31"
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("D40").Value = "1.74"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "38.89"
$ws.Range("E42").Value = "  -3.88%  "
$ws.Range("D43").Value = "151.37"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "20.54"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "0.599"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "1.76"
$ws.Range("E51").Value = "  +0.23%  "
